$d3text = @"
Background
The public health response to Zika outbreak has mostly focused on epidemiological surveillance, vector control, and individual level preventative measures.
 This qualitative study employs a social-ecological framework to examine how macro (historical, legislative, political, socio-economic factors), meso (sources of information, social support, social mobilization) and micro level factors (individual actions, behavioral changes) interacted to influence the response and behavior of women with respect to Zika in different contexts.
Methods
A qualitative study was carried out.
 Women were recruited through the snowball sampling technique from various locations in Brazil, Puerto Rico, and the United States.
 They were of different nationalities and ethnicities.
 Data were collected through semi-structured interviews.
 The data transcripts were analyzed using thematic analysis.
Results
Women in this study deemed the information provided as insufficient, which led them to actively reach out and access a variety of media sources.
 Social networks played a vital role in sharing information but also resulted in the spread of hoaxes or rumors.
 Participants in our research perceived socio-economic inequities but focused on how to remedy their microenvironments.
 They did not engage in major social activities.
 Lack of trust in governments placed women in vulnerable situations by preventing them to follow the guidance of health authorities.
 These impacts were also a result of the response tactics of health and government administrations in their failed attempts to ensure the well-being of their countries’ populations.
Conclusions
Our findings call for public health interventions that go beyond individual level behavioral change campaigns, to more comprehensively address the broader meso and macro level factors that influence women’ willingness and possibility to protect themselves.

"@

$e3text = @"
[Ana Rosa%Linde-Arias%linde14@yahoo.com%1,      Maria%Roura%NULL%1,      Eduardo%Siqueira%NULL%1]
"@

$d4text = @"
Background
Zika virus (ZIKV) infection during pregnancy has severe consequences on the new-born.
 The World Health Organization declared the Zika outbreak to be a Public Health Emergency of International Concern (PHEIC) in 2016. Health facilities in the regions most affected by Zika lacked the capacity to respond to the increased demand for contraception.
 The objectives were to explore healthcare users’ perceptions regarding contraception, Zika prevention during pregnancy and post-abortion care (PAC) services in the context of a Zika outbreak in Tegucigalpa, Honduras, and to follow these services over time.
Methods
This study was part of a broader implementation research study.
 We used qualitative research consistent with grounded theory approach.
 Semi-structured interviews and focus groups were performed with women and their partners who used contraceptive services or received PAC services.
 Data were collected in two stages from December 2017 to July 2018. Themes explored included contraception, Zika and PAC services.
Results
Participants had positive attitude towards the use of contraceptive methods and demanded more information on safety, efficacy and on side effects.
 Health care services were inconsistent in the provision of information on Zika and contraception services.
 ZIKV vector transmission was known but fewer participants were aware of risk of sexual transmission of Zika.
 Barriers to access healthcare services included contraceptive and PAC services included distance to healthcare facilities, disorganized admission process, long waiting times and out-of-pocket expenditure to purchase medicines.
 Furthermore, poor quality, mistreatment and abuse of women seeking PAC was prevalent.
 Some positive changes were noted over time, such as improvements in infrastructure including improved privacy and cleanliness, removal of fees, requisite to bring clean water to hospital.
Conclusions
Our results highlight the challenges and areas for improvement in policy and practice related to contraceptive services and PAC in the context of ZIKV infection.
 Public policies to prevent epidemics should focus more on providing proper sanitation; removing barriers to access and use of effective contraception as human rights priority.
 Zika epidemic has highlighted weaknesses in health systems that obstruct access to and use of sexual and reproductive health services.

"@

$e4text = @"
[Maria%Belizan%mbelizan@iecs.org.ar%1,      Edna%Maradiaga%edjamar3006@yahoo.com%1,      Javier%Roberti%jroberti@iecs.org.ar%1,      Maricela%Casco-Aguilar%marykasco@yahoo.com%1,      Alison F.%Ortez%alison_fabiola@yahoo.es%1,      Juan C.%Avila-Flores%javilaflores3@gmail.com%1,      Gloria%González%marilyntoin@yahoo.com%1,      Carolina%Bustillo%mcbu1502@yahoo.com%1,      Alejandra%Calderón%lilianalecalderon@gmail.com%1,      Harry%Bock%hbockme@hotmail.com%1,      María L.%Cafferata%NULL%1,      Adriano B.%Tavares%adriano.b.tavares@gmail.com%1,      Jackeline%Alger%jackelinealger@gmail.com%1,      Moazzam%Ali%alimoa@who.int%1]
"@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = $d3text
$ws.Range("E3").Value = $e3text
$ws.Range("D4").Value = $d4text
$ws.Range("E4").Value = $e4text
